# Add a new "AuthCapture" worksheet (a copy of "Auth") capturing a newer
# VLink smoke-test run, and update the selection/active-tab state to match.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Auth")

# Duplicate the Auth sheet (keeps formatting/columns/styles/data identical)
# and place the copy immediately after "Auth".
[void]$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item($ws1.Index + 1)
$ws2.Name = "AuthCapture"

# Update the "Date" column (B) on the new sheet with the timestamps from the
# newer VLink AuthCapture smoke-test execution.
$ws2.Range("B2").Value = "Tue Oct 04 21:04:53 EDT 2022"
$ws2.Range("B3").Value = "Tue Oct 04 21:05:17 EDT 2022"
$ws2.Range("B4").Value = "Tue Oct 04 21:05:40 EDT 2022"
$ws2.Range("B5").Value = "Tue Oct 04 21:06:02 EDT 2022"
$ws2.Range("B6").Value = "Tue Oct 04 21:06:25 EDT 2022"

# Restore the old sheet's selection to "select all columns" (as left by the
# Katalon run) and make the new AuthCapture sheet the active/selected tab.
[void]$ws1.Activate()
[void]$ws1.Range("A1:XFD1048576").Select()

[void]$ws2.Activate()
[void]$ws2.Range("C11").Select()
